# Update LR-pair data table (Pomc-Mc5r) following Dr Hou advice.
# Existing rows 2-7 get revised values (ligand-expressing/receptor-expressing
# cell counts and derived stats), and new target-cluster rows for "M2" and
# "sCs" clusters are added, extending the table down to row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Pomc"
$ws.Cells.Item(2, 3).Value = "Mc5r"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.4951615
$ws.Cells.Item(2, 8).Value = 2.990323
$ws.Cells.Item(2, 9).Value = 0.6020739711267923
$ws.Cells.Item(2, 10).Value = 0.5021622551131893
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.9666250000000001
$ws.Cells.Item(2, 14).Value = 1.93325
$ws.Cells.Item(2, 15).Value = 0.322724746715536
$ws.Cells.Item(2, 16).Value = 0.2550922792114143
$ws.Cells.Item(2, 17).Value = 1.4452604849375
$ws.Cells.Item(2, 18).Value = 5.781041939750001
$ws.Cells.Item(2, 19).Value = 0.194304169835911
$ws.Cells.Item(2, 20).Value = 0.1280977141907671

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pomc"
$ws.Cells.Item(3, 3).Value = "Mc5r"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.4951615
$ws.Cells.Item(3, 8).Value = 2.990323
$ws.Cells.Item(3, 9).Value = 0.6020739711267923
$ws.Cells.Item(3, 10).Value = 0.5021622551131893
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.093967666666667
$ws.Cells.Item(3, 14).Value = 3.281903
$ws.Cells.Item(3, 15).Value = 0.3652403342971534
$ws.Cells.Item(3, 16).Value = 0.4330470018987601
$ws.Cells.Item(3, 17).Value = 1.635658337444833
$ws.Cells.Item(3, 18).Value = 9.813950024669
$ws.Cells.Item(3, 19).Value = 0.2199016984859643
$ws.Cells.Item(3, 20).Value = 0.2174598590434869

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Pomc"
$ws.Cells.Item(4, 3).Value = "Mc5r"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.4951615
$ws.Cells.Item(4, 8).Value = 2.990323
$ws.Cells.Item(4, 9).Value = 0.6020739711267923
$ws.Cells.Item(4, 10).Value = 0.5021622551131893
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.01150933333333333
$ws.Cells.Item(4, 14).Value = 0.034528
$ws.Cells.Item(4, 15).Value = 0.003842593234051133
$ws.Cells.Item(4, 16).Value = 0.004555968558961185
$ws.Cells.Item(4, 17).Value = 0.01720831209066667
$ws.Cells.Item(4, 18).Value = 0.103249872544
$ws.Cells.Item(4, 19).Value = 0.002313525367850109
$ws.Cells.Item(4, 20).Value = 0.002287835445792736

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Pomc"
$ws.Cells.Item(5, 3).Value = "Mc5r"
$ws.Cells.Item(5, 4).Value = "Neutro"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.4951615
$ws.Cells.Item(5, 8).Value = 2.990323
$ws.Cells.Item(5, 9).Value = 0.6020739711267923
$ws.Cells.Item(5, 10).Value = 0.5021622551131893
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.482754
$ws.Cells.Item(5, 14).Value = 1.448262
$ws.Cells.Item(5, 15).Value = 0.1611759083159569
$ws.Cells.Item(5, 16).Value = 0.1910981272340779
$ws.Cells.Item(5, 17).Value = 0.721795194771
$ws.Cells.Item(5, 18).Value = 4.330771168626
$ws.Cells.Item(5, 19).Value = 0.09703981916975597
$ws.Cells.Item(5, 20).Value = 0.09596226651977176

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Pomc"
$ws.Cells.Item(6, 3).Value = "Mc5r"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.4951615
$ws.Cells.Item(6, 8).Value = 2.990323
$ws.Cells.Item(6, 9).Value = 0.6020739711267923
$ws.Cells.Item(6, 10).Value = 0.5021622551131893
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.5
$ws.Cells.Item(6, 13).Value = 0.4403435
$ws.Cells.Item(6, 14).Value = 0.880687
$ws.Cells.Item(6, 15).Value = 0.1470164174373026
$ws.Cells.Item(6, 16).Value = 0.1162066230967866
$ws.Cells.Item(6, 17).Value = 0.65838464797525
$ws.Cells.Item(6, 18).Value = 2.633538591901
$ws.Cells.Item(6, 19).Value = 0.08851475826731098
$ws.Cells.Item(6, 20).Value = 0.05835457991337081

# Row 7
$ws.Cells.Item(7, 1).Value = "Neutro"
$ws.Cells.Item(7, 2).Value = "Pomc"
$ws.Cells.Item(7, 3).Value = "Mc5r"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.9881903333333333
$ws.Cells.Item(7, 8).Value = 2.964571
$ws.Cells.Item(7, 9).Value = 0.3979260288732077
$ws.Cells.Item(7, 10).Value = 0.4978377448868108
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.9666250000000001
$ws.Cells.Item(7, 14).Value = 1.93325
$ws.Cells.Item(7, 15).Value = 0.322724746715536
$ws.Cells.Item(7, 16).Value = 0.2550922792114143
$ws.Cells.Item(7, 17).Value = 0.9552094809583334
$ws.Cells.Item(7, 18).Value = 5.73125688575
$ws.Cells.Item(7, 19).Value = 0.128420576879625
$ws.Cells.Item(7, 20).Value = 0.1269945650206472

# Row 8
$ws.Cells.Item(8, 1).Value = "Neutro"
$ws.Cells.Item(8, 2).Value = "Pomc"
$ws.Cells.Item(8, 3).Value = "Mc5r"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.9881903333333333
$ws.Cells.Item(8, 8).Value = 2.964571
$ws.Cells.Item(8, 9).Value = 0.3979260288732077
$ws.Cells.Item(8, 10).Value = 0.4978377448868108
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.093967666666667
$ws.Cells.Item(8, 14).Value = 3.281903
$ws.Cells.Item(8, 15).Value = 0.3652403342971534
$ws.Cells.Item(8, 16).Value = 0.4330470018987601
$ws.Cells.Item(8, 17).Value = 1.081048273179222
$ws.Cells.Item(8, 18).Value = 9.729434458613
$ws.Cells.Item(8, 19).Value = 0.1453386358111891
$ws.Cells.Item(8, 20).Value = 0.2155871428552732

# Row 9
$ws.Cells.Item(9, 1).Value = "Neutro"
$ws.Cells.Item(9, 2).Value = "Pomc"
$ws.Cells.Item(9, 3).Value = "Mc5r"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.9881903333333333
$ws.Cells.Item(9, 8).Value = 2.964571
$ws.Cells.Item(9, 9).Value = 0.3979260288732077
$ws.Cells.Item(9, 10).Value = 0.4978377448868108
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.01150933333333333
$ws.Cells.Item(9, 14).Value = 0.034528
$ws.Cells.Item(9, 15).Value = 0.003842593234051133
$ws.Cells.Item(9, 16).Value = 0.004555968558961185
$ws.Cells.Item(9, 17).Value = 0.01137341194311111
$ws.Cells.Item(9, 18).Value = 0.102360707488
$ws.Cells.Item(9, 19).Value = 0.001529067866201024
$ws.Cells.Item(9, 20).Value = 0.002268133113168449

# Row 10
$ws.Cells.Item(10, 1).Value = "Neutro"
$ws.Cells.Item(10, 2).Value = "Pomc"
$ws.Cells.Item(10, 3).Value = "Mc5r"
$ws.Cells.Item(10, 4).Value = "Neutro"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.9881903333333333
$ws.Cells.Item(10, 8).Value = 2.964571
$ws.Cells.Item(10, 9).Value = 0.3979260288732077
$ws.Cells.Item(10, 10).Value = 0.4978377448868108
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.482754
$ws.Cells.Item(10, 14).Value = 1.448262
$ws.Cells.Item(10, 15).Value = 0.1611759083159569
$ws.Cells.Item(10, 16).Value = 0.1910981272340779
$ws.Cells.Item(10, 17).Value = 0.4770528361779999
$ws.Cells.Item(10, 18).Value = 4.293475525601999
$ws.Cells.Item(10, 19).Value = 0.06413608914620096
$ws.Cells.Item(10, 20).Value = 0.09513586071430619

# Row 11
$ws.Cells.Item(11, 1).Value = "Neutro"
$ws.Cells.Item(11, 2).Value = "Pomc"
$ws.Cells.Item(11, 3).Value = "Mc5r"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.9881903333333333
$ws.Cells.Item(11, 8).Value = 2.964571
$ws.Cells.Item(11, 9).Value = 0.3979260288732077
$ws.Cells.Item(11, 10).Value = 0.4978377448868108
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.5
$ws.Cells.Item(11, 13).Value = 0.4403435
$ws.Cells.Item(11, 14).Value = 0.880687
$ws.Cells.Item(11, 15).Value = 0.1470164174373026
$ws.Cells.Item(11, 16).Value = 0.1162066230967866
$ws.Cells.Item(11, 17).Value = 0.4351431900461666
$ws.Cells.Item(11, 18).Value = 2.610859140277
$ws.Cells.Item(11, 19).Value = 0.05850165916999164
$ws.Cells.Item(11, 20).Value = 0.05785204318341583
